# Update cryptocurrency price/volume data (and swap Filecoin/Maker rows 43-44)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.850.04"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.515.16"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D5").Value = "'600.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").Value = "'181.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.51%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "3.513.89"

$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("E10").Value = "  +6.56%  "

$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("D12").Value = "'0.438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").Value = "4.123.78"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'32.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.73%  "

$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").Value = "67.886.08"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "'0.0000180"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "3.526.98"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").Value = "'14.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").Value = "'400.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").Value = "'7.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").Value = "'73.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").Value = "'0.543"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("D28").Value = "'10.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "

$ws.Range("E29").Value = "  -2.39%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "'6.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("D33").Value = "'2.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").Value = "'23.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("E35").Value = "  +1.65%  "

$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("D38").Value = "'163.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("D39").Value = "'0.880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("D41").Value = "'2.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.62%  "

$ws.Range("D42").Value = "'6.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.43%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.890.22"
$ws.Range("E43").Value = "  +2.71%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'0.0736"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "

$ws.Range("D46").Value = "'26.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").Value = "'26.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").Value = "'42.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("D49").Value = "'347.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").Value = "'0.0304"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("E51").Value = "  -1.40%  "
